$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header text in A1 (shared string content change)
$ws.Range("A1").Value = "Glucosio.totale"

# The trailing empty/styled row 20 is no longer part of the data -- remove it
# so the sheet's used range shrinks back to A1:A19.
$ws.Range("A20").EntireRow.Delete()
